$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing "Open Sans 10pt / FF333333" cell style (style index 1)
# by copying formats from an already-styled cell instead of setting Font
# properties one-by-one (which would create new transient style/font
# records even though the final cell ends up matching style 1).
$fmtSource = $ws.Range("C11")

# --- New account-area rows -------------------------------------------------
# Values are written in the same order the author typed them so the shared
# string table grows in the same sequence (ac_button_delete_account,
# areaAccount, ac_email, ac_verified).
$ws.Range("C20").Value2 = "ac_button_delete_account"
$fmtSource.Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("B20").Value2 = "areaAccount"

$ws.Range("C22").Value2 = "ac_email"
$fmtSource.Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("C21").Value2 = "ac_verified"
$fmtSource.Copy()
$ws.Range("C21").PasteSpecial(-4122)

# B21 used to hold "areaMessages" - clear the text but keep the styled
# (Open Sans) formatting that was already there.
$ws.Range("B21").ClearContents()
$fmtSource.Copy()
$ws.Range("B21").PasteSpecial(-4122)

# New trailing styled-but-empty row.
$ws.Range("B23").Value2 = ""
$fmtSource.Copy()
$ws.Range("B23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Match the saved selection / scroll position.
$ws.Range("C19").Select()
